$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header for left table: "interesting" -> "negative"
$ws.Range("A1").Value = "negative"

# --- Rows 3-7 (left table columns B,C,D,H): refreshed counts from the larger dataset ---
$ws.Range("B3").Value = 0.8529411764705882
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 29
$ws.Range("H3").Value = 5
$ws.Range("B4").Value = 0.6061643835616438
$ws.Range("C4").Value = 177
$ws.Range("D4").Value = 177
$ws.Range("H4").Value = 115
$ws.Range("B5").Value = 0.2151162790697674
$ws.Range("C5").Value = 111
$ws.Range("D5").Value = 111
$ws.Range("H5").Value = 405
$ws.Range("B6").Value = 0.1957671957671958
$ws.Range("C6").Value = 37
$ws.Range("D6").Value = 37
$ws.Range("H6").Value = 152
$ws.Range("B7").Value = 0.07222222222222222
$ws.Range("C7").Value = 26
$ws.Range("D7").Value = 26
$ws.Range("H7").Value = 334

# --- Row 8 no longer has a left-table entry; clear A8:H8 entirely (style + value) ---
$ws.Range("A8:H8").Clear()

# --- Extend right table (J:Q) formatting down through row 45 (copy format from row 29) ---
$ws.Range("J29:Q29").Copy()
$ws.Range("J30:Q45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Right table (J:Q) rows 3-45: word, +%, type occurences, total occurences, both%, -%, both flag, total "both" occurences ---
$ws.Range("J3").Value = "interesting"
$ws.Range("K3").Value = 0.9696969696969697
$ws.Range("L3").Value = 32
$ws.Range("M3").Value = 32
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1
$ws.Range("J4").Value = "happy"
$ws.Range("K4").Value = 0.9615384615384616
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 1
$ws.Range("J5").Value = "best"
$ws.Range("K5").Value = 0.9322033898305084
$ws.Range("L5").Value = 55
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 4
$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.875
$ws.Range("L6").Value = 98
$ws.Range("M6").Value = 98
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 14
$ws.Range("J7").Value = "love"
$ws.Range("K7").Value = 0.8695652173913043
$ws.Range("L7").Value = 40
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 7
$ws.Range("J8").Value = "won"
$ws.Range("K8").Value = 0.8205128205128205
$ws.Range("L8").Value = 32
$ws.Range("M8").Value = 32
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 7
$ws.Range("J9").Value = "positive"
$ws.Range("K9").Value = 0.8103448275862069
$ws.Range("L9").Value = 47
$ws.Range("M9").Value = 47
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 11
$ws.Range("J10").Value = "thanks"
$ws.Range("K10").Value = 0.8048780487804879
$ws.Range("L10").Value = 66
$ws.Range("M10").Value = 66
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 16
$ws.Range("J11").Value = "thank"
$ws.Range("K11").Value = 0.796875
$ws.Range("L11").Value = 102
$ws.Range("M11").Value = 102
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 26
$ws.Range("J12").Value = "confidence"
$ws.Range("K12").Value = 0.75
$ws.Range("L12").Value = 27
$ws.Range("M12").Value = 27
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 9
$ws.Range("J13").Value = "special"
$ws.Range("K13").Value = 0.75
$ws.Range("L13").Value = 27
$ws.Range("M13").Value = 27
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 9
$ws.Range("J14").Value = "free"
$ws.Range("K14").Value = 0.75
$ws.Range("L14").Value = 90
$ws.Range("M14").Value = 90
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 30
$ws.Range("J15").Value = "safe"
$ws.Range("K15").Value = 0.7183098591549296
$ws.Range("L15").Value = 102
$ws.Range("M15").Value = 102
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 40
$ws.Range("J16").Value = "support"
$ws.Range("K16").Value = 0.7169811320754716
$ws.Range("L16").Value = 76
$ws.Range("M16").Value = 76
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 30
$ws.Range("J17").Value = "safety"
$ws.Range("K17").Value = 0.7058823529411765
$ws.Range("L17").Value = 36
$ws.Range("M17").Value = 36
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 15
$ws.Range("J18").Value = "good"
$ws.Range("K18").Value = 0.7
$ws.Range("L18").Value = 112
$ws.Range("M18").Value = 112
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 48
$ws.Range("J19").Value = "better"
$ws.Range("K19").Value = 0.6666666666666666
$ws.Range("L19").Value = 42
$ws.Range("M19").Value = 42
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 21
$ws.Range("J20").Value = "relief"
$ws.Range("K20").Value = 0.6
$ws.Range("L20").Value = 30
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 20
$ws.Range("J21").Value = "well"
$ws.Range("K21").Value = 0.5638297872340425
$ws.Range("L21").Value = 53
$ws.Range("M21").Value = 53
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 41
$ws.Range("J22").Value = "heroes"
$ws.Range("K22").Value = 0.5531914893617021
$ws.Range("L22").Value = 26
$ws.Range("M22").Value = 26
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 21
$ws.Range("J23").Value = "hand"
$ws.Range("K23").Value = 0.5430809399477807
$ws.Range("L23").Value = 208
$ws.Range("M23").Value = 208
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 175
$ws.Range("J24").Value = "fresh"
$ws.Range("K24").Value = 0.5416666666666666
$ws.Range("L24").Value = 26
$ws.Range("M24").Value = 26
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 22
$ws.Range("J25").Value = "care"
$ws.Range("K25").Value = 0.4606741573033708
$ws.Range("L25").Value = 41
$ws.Range("M25").Value = 41
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 48
$ws.Range("J26").Value = "like"
$ws.Range("K26").Value = 0.4588235294117647
$ws.Range("L26").Value = 156
$ws.Range("M26").Value = 156
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 184
$ws.Range("J27").Value = "help"
$ws.Range("K27").Value = 0.4542372881355932
$ws.Range("L27").Value = 134
$ws.Range("M27").Value = 134
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 161
$ws.Range("J28").Value = "increase"
$ws.Range("K28").Value = 0.3974358974358974
$ws.Range("L28").Value = 31
$ws.Range("M28").Value = 31
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 47
$ws.Range("J29").Value = "please"
$ws.Range("K29").Value = 0.3849372384937239
$ws.Range("L29").Value = 92
$ws.Range("M29").Value = 92
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 147
$ws.Range("J30").Value = "hope"
$ws.Range("K30").Value = 0.3846153846153846
$ws.Range("L30").Value = 25
$ws.Range("M30").Value = 25
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 40
$ws.Range("J31").Value = "protect"
$ws.Range("K31").Value = 0.3698630136986301
$ws.Range("L31").Value = 27
$ws.Range("M31").Value = 27
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = $false
$ws.Range("Q31").Value = 46
$ws.Range("J32").Value = "store"
$ws.Range("K32").Value = 0.03691275167785235
$ws.Range("L32").Value = 33
$ws.Range("M32").Value = 33
$ws.Range("N32").Value = 1
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = $false
$ws.Range("Q32").Value = 861
$ws.Range("J33").Value = "you"
$ws.Range("K33").Value = 0.03333333333333333
$ws.Range("L33").Value = 40
$ws.Range("M33").Value = 40
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = $false
$ws.Range("Q33").Value = 1160
$ws.Range("J34").Value = "!"
$ws.Range("K34").Value = 0.03140830800405268
$ws.Range("L34").Value = 31
$ws.Range("M34").Value = 33
$ws.Range("N34").Value = 0.94
$ws.Range("O34").Value = 0.06000000000000005
$ws.Range("P34").Value = $true
$ws.Range("Q34").Value = 956
$ws.Range("J35").Value = "grocery"
$ws.Range("K35").Value = 0.02996670366259711
$ws.Range("L35").Value = 27
$ws.Range("M35").Value = 27
$ws.Range("N35").Value = 1
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = $false
$ws.Range("Q35").Value = 874
$ws.Range("J36").Value = "for"
$ws.Range("K36").Value = 0.02199874292897549
$ws.Range("L36").Value = 35
$ws.Range("M36").Value = 38
$ws.Range("N36").Value = 0.92
$ws.Range("O36").Value = 0.07999999999999996
$ws.Range("P36").Value = $true
$ws.Range("Q36").Value = 1556
$ws.Range("J37").Value = "and"
$ws.Range("K37").Value = 0.02134032197678772
$ws.Range("L37").Value = 57
$ws.Range("M37").Value = 59
$ws.Range("N37").Value = 0.97
$ws.Range("O37").Value = 0.03000000000000003
$ws.Range("P37").Value = $true
$ws.Range("Q37").Value = 2614
$ws.Range("J38").Value = ","
$ws.Range("K38").Value = 0.0204582651391162
$ws.Range("L38").Value = 50
$ws.Range("M38").Value = 51
$ws.Range("N38").Value = 0.98
$ws.Range("O38").Value = 0.02000000000000002
$ws.Range("P38").Value = $true
$ws.Range("Q38").Value = 2394
$ws.Range("J39").Formula = "=""19"""
$ws.Range("J39").Copy()
$ws.Range("J39").PasteSpecial(-4163)
$ws.Range("K39").Value = 0.01637043966323667
$ws.Range("L39").Value = 35
$ws.Range("M39").Value = 40
$ws.Range("N39").Value = 0.88
$ws.Range("O39").Value = 0.12
$ws.Range("P39").Value = $true
$ws.Range("Q39").Value = 2103
$ws.Range("J40").Value = "."
$ws.Range("K40").Value = 0.01543086172344689
$ws.Range("L40").Value = 77
$ws.Range("M40").Value = 83
$ws.Range("N40").Value = 0.93
$ws.Range("O40").Value = 0.06999999999999995
$ws.Range("P40").Value = $true
$ws.Range("Q40").Value = 4913
$ws.Range("J41").Value = "a"
$ws.Range("K41").Value = 0.01387082791504118
$ws.Range("L41").Value = 32
$ws.Range("M41").Value = 34
$ws.Range("N41").Value = 0.94
$ws.Range("O41").Value = 0.06000000000000005
$ws.Range("P41").Value = $true
$ws.Range("Q41").Value = 2275
$ws.Range("J42").Value = "of"
$ws.Range("K42").Value = 0.01306363253265908
$ws.Range("L42").Value = 31
$ws.Range("M42").Value = 37
$ws.Range("N42").Value = 0.84
$ws.Range("O42").Value = 0.16
$ws.Range("P42").Value = $true
$ws.Range("Q42").Value = 2342
$ws.Range("J43").Value = "to"
$ws.Range("K43").Value = 0.01294199214236191
$ws.Range("L43").Value = 56
$ws.Range("M43").Value = 57
$ws.Range("N43").Value = 0.98
$ws.Range("O43").Value = 0.02000000000000002
$ws.Range("P43").Value = $true
$ws.Range("Q43").Value = 4271
$ws.Range("J44").Value = "the"
$ws.Range("K44").Value = 0.01163016088389223
$ws.Range("L44").Value = 60
$ws.Range("M44").Value = 66
$ws.Range("N44").Value = 0.91
$ws.Range("O44").Value = 0.08999999999999997
$ws.Range("P44").Value = $true
$ws.Range("Q44").Value = 5099
$ws.Range("J45").Value = "co"
$ws.Range("K45").Value = 0.009990331936835321
$ws.Range("L45").Value = 31
$ws.Range("M45").Value = 35
$ws.Range("N45").Value = 0.89
$ws.Range("O45").Value = 0.11
$ws.Range("P45").Value = $true
$ws.Range("Q45").Value = 3072
$excel.CutCopyMode = $false
